$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1208
$ws.Range("I19").Value = 1177.4445
$ws.Range("J19").Value = 1299.6666
$ws.Range("K19").Value = 1177.4445
$ws.Range("L19").Value = 1299.6666
$ws.Range("M19").Value = -1002.4445
$ws.Range("N19").Value = -1649.6666
$ws.Range("H88").Value = 1665
$ws.Range("J88").Value = 1665
$ws.Range("L88").Value = 1665
$ws.Range("N88").Value = -2477
$ws.Range("H91").Value = 1665
$ws.Range("J91").Value = 1665
$ws.Range("L91").Value = 1665
$ws.Range("N91").Value = -4473
$ws.Range("H138").Value = 2175.25
$ws.Range("I138").Value = 2175.25
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 6525.75
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -1385.75
$ws.Range("N138").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1358.3334
$ws.Range("I2").Value = 1358.3334
$ws.Range("K2").Value = 1358.3334
$ws.Range("M2").Value = -1245.3334
$ws.Range("H74").Value = 2565.5557
$ws.Range("I74").Value = 1482.6666
$ws.Range("J74").Value = 4731.3335
$ws.Range("K74").Value = 1482.6666
$ws.Range("L74").Value = 4731.3335
$ws.Range("M74").Value = -608.6666
$ws.Range("N74").Value = -6479.3335
$ws.Range("H77").Value = 2565.5557
$ws.Range("I77").Value = 1482.6666
$ws.Range("J77").Value = 4731.3335
$ws.Range("K77").Value = 7413.333000000001
$ws.Range("L77").Value = 23656.6675
$ws.Range("M77").Value = -3045.333000000001
$ws.Range("N77").Value = -32392.6675
$ws.Range("H116").Value = 1358.3334
$ws.Range("I116").Value = 1358.3334
$ws.Range("K116").Value = 1358.3334
$ws.Range("M116").Value = 935.6666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1358.3334
$ws.Range("I3").Value = 1358.3334
$ws.Range("K3").Value = 1358.3334
$ws.Range("M3").Value = -1244.3334
$ws.Range("H20").Value = 5158.4
$ws.Range("I20").Value = 5158.4
$ws.Range("K20").Value = 5158.4
$ws.Range("M20").Value = -4911.4
$ws.Range("H94").Value = 1183.6364
$ws.Range("I94").Value = 1202
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1202
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -751
$ws.Range("N94").Value = -1902
$ws.Range("H99").Value = 1342.5
$ws.Range("I99").Value = 1342.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1342.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 155.5
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 29568.857
$ws.Range("I105").Value = 1196.8
$ws.Range("J105").Value = 100499
$ws.Range("K105").Value = 1196.8
$ws.Range("L105").Value = 100499
$ws.Range("M105").Value = 550.2
$ws.Range("N105").Value = -103993
$ws.Range("H134").Value = 2597
$ws.Range("I134").Value = 2597
$ws.Range("K134").Value = 7791
$ws.Range("M134").Value = -5256

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 131.5
$ws.Range("I7").Value = 170.33333
$ws.Range("K7").Value = 170.33333
$ws.Range("M7").Value = -57.33332999999999
$ws.Range("H58").Value = 1815.4546
$ws.Range("I58").Value = 807.125
$ws.Range("K58").Value = 807.125
$ws.Range("M58").Value = -604.125
$ws.Range("H99").Value = 6571.1
$ws.Range("I99").Value = 6666.6665
$ws.Range("J99").Value = 6530.143
$ws.Range("K99").Value = 6666.6665
$ws.Range("L99").Value = 6530.143
$ws.Range("M99").Value = -5168.6665
$ws.Range("N99").Value = -9526.143
$ws.Range("H123").Value = 150999
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 6571.1
$ws.Range("I126").Value = 6666.6665
$ws.Range("J126").Value = 6530.143
$ws.Range("K126").Value = 19999.9995
$ws.Range("L126").Value = 19590.429
$ws.Range("M126").Value = -17529.9995
$ws.Range("N126").Value = -24530.429
$ws.Range("H132").Value = 2794.077
$ws.Range("I132").Value = 2976.6667
$ws.Range("J132").Value = 2383.25
$ws.Range("K132").Value = 8930.000100000001
$ws.Range("L132").Value = 7149.75
$ws.Range("M132").Value = -6400.000100000001
$ws.Range("N132").Value = -12209.75
$ws.Range("H134").Value = 743.25
$ws.Range("I134").Value = 735.1429000000001
$ws.Range("K134").Value = 2205.4287
$ws.Range("M134").Value = 329.5712999999996
$ws.Range("H136").Value = 1815.4546
$ws.Range("I136").Value = 807.125
$ws.Range("K136").Value = 2421.375
$ws.Range("M136").Value = 128.625

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2299.8
$ws.Range("J39").Value = 2499.75
$ws.Range("L39").Value = 7499.25
$ws.Range("N39").Value = -8087.25
$ws.Range("H55").Value = 576
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("H58").Value = 6250
$ws.Range("I58").Value = 6250
$ws.Range("K58").Value = 18750
$ws.Range("M58").Value = -18622
$ws.Range("H113").Value = 751.1667
$ws.Range("I113").Value = 614.25
$ws.Range("J113").Value = 1025
$ws.Range("K113").Value = 1842.75
$ws.Range("L113").Value = 3075
$ws.Range("M113").Value = 327.25
$ws.Range("N113").Value = -7415
$ws.Range("H114").Value = 1200.5
$ws.Range("H118").Value = 3300
$ws.Range("I118").Value = 3300
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 9900
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -8657
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 549.5
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 8991
$ws.Range("N122").Value = -13891
$ws.Range("H140").Value = 3666.3333
$ws.Range("I140").Value = 2999.5
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 8998.5
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -3818.5
$ws.Range("N140").Value = -25360

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2750
$ws.Range("I7").Value = 2750
$ws.Range("K7").Value = 2750
$ws.Range("M7").Value = -2638
$ws.Range("H40").Value = 2999.5
$ws.Range("J40").Value = 3999
$ws.Range("L40").Value = 3999
$ws.Range("N40").Value = -4271
$ws.Range("H55").Value = 907.3333
$ws.Range("I55").Value = 483.375
$ws.Range("J55").Value = 1391.8572
$ws.Range("K55").Value = 483.375
$ws.Range("L55").Value = 1391.8572
$ws.Range("M55").Value = -310.375
$ws.Range("N55").Value = -1737.8572
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780
$ws.Range("H136").Value = 5888
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5888
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 17664
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -22764

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 8600
$ws.Range("J50").Value = 8600
$ws.Range("L50").Value = 8600
$ws.Range("N50").Value = -9862
$ws.Range("H113").Value = 4464.1
$ws.Range("J113").Value = 4731.1113
$ws.Range("L113").Value = 14193.3339
$ws.Range("N113").Value = -18533.3339
$ws.Range("H136").Value = 7447
$ws.Range("I136").Value = 7447
$ws.Range("K136").Value = 22341
$ws.Range("M136").Value = -19791
